$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.003863092951775041
$ws.Range("B2").Value = -0.03219275312800293

$ws.Range("A3").Value = 0.06728693851364928
$ws.Range("B3").Value = 0.07015337425915152

$ws.Range("A4").Value = 0.0009801596384573432
$ws.Range("B4").Value = -0.03002888080940301

$ws.Range("A5").Value = 0.1691303632187597
$ws.Range("B5").Value = 0.1747496709315947

$ws.Range("A6").Value = -0.0636064919738203
$ws.Range("B6").Value = -0.03191399942666489

$ws.Range("A7").Value = 0.5014689555692808
$ws.Range("B7").Value = 0.4913937796165608

$ws.Range("A8").Value = 0.2572704901080778
$ws.Range("B8").Value = 0.263703175192693
